$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.454.88"
$ws.Range("E2").Value = "  +2.38%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.99"
$ws.Range("E3").Value = "  +2.05%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.78"
$ws.Range("E5").Value = "  +0.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.610"
$ws.Range("E6").Value = "  +1.87%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.07"
$ws.Range("E8").Value = "  +8.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.307"
$ws.Range("E9").Value = "  +5.23%  "

$ws.Range("E10").Value = "  +1.10%  "

$ws.Range("E11").Value = "  +3.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.114.01"
$ws.Range("E12").Value = "  +1.90%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.59"
$ws.Range("E13").Value = "  +2.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.852.15"
$ws.Range("E14").Value = "  +1.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.72"
$ws.Range("E15").Value = "  +6.37%  "

$ws.Range("E16").Value = "  +4.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.383.49"
$ws.Range("E17").Value = "  +2.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.98"
$ws.Range("E18").Value = "  +1.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.99"
$ws.Range("E19").Value = "  +0.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0792"
$ws.Range("E20").Value = "  +1.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.14"
$ws.Range("E21").Value = "  +7.81%  "

$ws.Range("E22").Value = "  +15.90%  "

$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.20"
$ws.Range("E24").Value = "  -0.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.98"
$ws.Range("E25").Value = "  -0.62%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.90"
$ws.Range("E26").Value = "  -0.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.87"
$ws.Range("E27").Value = "  +2.82%  "

$ws.Range("E28").Value = "  +0.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.498.32"
$ws.Range("E29").Value = "  +43.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.01"
$ws.Range("E30").Value = "  +0.46%  "

$ws.Range("E31").Value = "  +7.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.93"
$ws.Range("E32").Value = "  +3.06%  "

$ws.Range("E33").Value = "  +2.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0533"
$ws.Range("E34").Value = "  +1.79%  "

$ws.Range("E35").Value = "  +2.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.674"
$ws.Range("E36").Value = "  +2.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "88.80"
$ws.Range("E37").Value = "  +9.50%  "

$ws.Range("E38").Value = "  +1.59%  "

$ws.Range("E39").Value = "  +9.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.336.79"
$ws.Range("E40").Value = "  -2.19%  "

$ws.Range("E41").Value = "  +3.22%  "

$ws.Range("E42").Value = "  +2.65%  "

$ws.Range("E43").Value = "  +4.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.90"
$ws.Range("E44").Value = "  +5.37%  "

$ws.Range("E45").Value = "  +1.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.84"
$ws.Range("E46").Value = "  +2.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0521"
$ws.Range("E47").Value = "  +3.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.012.76"
$ws.Range("E48").Value = "  +1.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.99"
$ws.Range("E49").Value = "  +2.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.20"
$ws.Range("E50").Value = "  +1.31%  "

$ws.Range("E51").Value = "  +0.09%  "
